# Update dashboards - 2025-11-14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style-only changes: highlight (fillId=3, style 48) vs plain (fillId=0, style 47) ---
# A cell with style 47 elsewhere (C3) is used as the "plain" format donor,
# and a cell with style 48 elsewhere (N47) is used as the "highlighted" format donor.

# Row 24 (Cons Credit - Revolving, PCEPI date C24): de-highlight -> style 47
$ws.Range("C3").Copy()
$ws.Range("C24").PasteSpecial(-4122)

# Row 25 (Cons Credit - NonRevolving, PCEPI date C25): de-highlight -> style 47
$ws.Range("C3").Copy()
$ws.Range("C25").PasteSpecial(-4122)

# Row 39 (Nominal Broad US Dollar Index date N39): de-highlight -> style 47
$ws.Range("C3").Copy()
$ws.Range("N39").PasteSpecial(-4122)

# Row 51 (30y Mtg. date N51): newly highlighted -> style 48
$ws.Range("N47").Copy()
$ws.Range("N51").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 28: Dur. Order M/M % Delta recomputed values ---
$ws.Range("F28").Value2 = 0.02933036907051223
$ws.Range("G28").Value2 = -0.02756598804853716

# --- Row 29: Dur. Order Y/Y % Delta (G) + 5yr,5yr Forward (T5YIFR) block ---
$ws.Range("G29").Value2 = 0.03387397095834726
$ws.Range("N29").Value2 = 45974
$ws.Range("T29").Value2 = $ws.Range("S29").Value2   # 2.2 -> T29
$ws.Range("S29").ClearContents()                    # S29 -> blank
$ws.Range("R29").Value2 = $ws.Range("Q29").Value2   # 2.17 -> R29
$ws.Range("Q29").Value2 = 2.19                       # new present value

# --- Row 30: 10yr TIPS (T10YIE) block ---
$ws.Range("N30").Value2 = 45974
$ws.Range("T30").Value2 = $ws.Range("S30").Value2   # 2.29 -> T30
$ws.Range("S30").ClearContents()                    # S30 -> blank
$ws.Range("R30").Value2 = $ws.Range("Q30").Value2   # 2.27 -> R30
$ws.Range("Q30").Value2 = 2.28                       # new present value

# --- Row 47: FFR date update only ---
$ws.Range("N47").Value2 = 45973

# --- Row 48: 2y UST (DGS2) block ---
$ws.Range("N48").Value2 = 45973
$ws.Range("S48").Value2 = $ws.Range("Q48").Value2   # 3.58 -> S48
$ws.Range("T48").ClearContents()                    # 3.55 dropped
$ws.Range("U48").ClearContents()                    # 3.57 dropped
$ws.Range("Q48").Value2 = 3.56                       # new present value

# --- Row 49: 5y UST (DGS5) block ---
$ws.Range("N49").Value2 = 45973
$ws.Range("S49").Value2 = $ws.Range("Q49").Value2   # 3.72 -> S49
$ws.Range("T49").ClearContents()                    # 3.67 dropped
$ws.Range("U49").ClearContents()                    # 3.69 dropped
$ws.Range("Q49").Value2 = 3.68                       # new present value

# --- Row 50: 10y UST (DGS10) block ---
$ws.Range("N50").Value2 = 45973
$ws.Range("S50").Value2 = $ws.Range("Q50").Value2   # 4.13 -> S50
$ws.Range("T50").ClearContents()                    # 4.11 dropped
$ws.Range("U50").ClearContents()                    # 4.11 dropped
$ws.Range("Q50").Value2 = 4.08                       # new present value

# --- Row 51: 30y Mtg. (MORTGAGE30US) block - shift right by one (weekly) ---
$ws.Range("U51").Value2 = $ws.Range("T51").Value2   # 6.27 -> U51
$ws.Range("T51").Value2 = $ws.Range("S51").Value2   # 6.19 -> T51
$ws.Range("S51").Value2 = $ws.Range("R51").Value2   # 6.17 -> S51
$ws.Range("R51").Value2 = $ws.Range("Q51").Value2   # 6.22 -> R51
$ws.Range("Q51").Value2 = 6.24                       # new present value
$ws.Range("N51").Value2 = 45971

# --- Row 52: BAA (DBAA) block ---
$ws.Range("N52").Value2 = 45973
$ws.Range("S52").Value2 = $ws.Range("Q52").Value2   # 5.86 -> S52
$ws.Range("T52").ClearContents()                    # 5.86 dropped
$ws.Range("U52").ClearContents()                    # 5.83 dropped
$ws.Range("Q52").Value2 = 5.83                       # new present value
